# Regenerate the "K" column (column G) values in the save_data sheet.
# These new values reflect a recalculation of strikeouts (K) replacing
# the previous "Strike#" derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 2
    6  = 1
    7  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 0
    13 = 2
    14 = 3
    15 = 1
    16 = 2
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 3
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 3
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 0
    46 = 2
    47 = 1
    48 = 2
    49 = 0
    50 = 1
    51 = 3
    52 = 1
    53 = 2
    54 = 1
    55 = 2
    56 = 1
    57 = 3
    58 = 1
    59 = 2
    60 = 0
    61 = 0
    62 = 1
    63 = 2
    64 = 1
    65 = 0
    66 = 4
    68 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
